$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking value per correct answer (row 11): 3 -> 5
$ws.Range("B11").Value = 5

# Update total correct score (row 12): 63 -> 105 (21 questions * 5 marks)
$ws.Range("B12").Value = 105

# Update the correct/total marks display text (row 12): "63/84" -> "105/140"
$ws.Range("E12").Value = "105/140"
